# Update the NATMI ligand-receptor pair sheet with newly-computed TPM values.
# The sheet originally held 4 result rows, all with "MuSCs" as the sending
# cluster. The new script run adds "ECs" as an additional sending cluster
# (4 more rows against the same 4 target clusters), and refreshes every
# expression/specificity metric for both sending clusters with the updated
# TPM-derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (after the A:D identifier columns) hold, in order:
#   E  Ligand-expressing cells
#   F  Ligand detection rate
#   G  Ligand average expression value
#   H  Ligand total expression value
#   I  Ligand derived specificity of average expression value
#   J  Ligand derived specificity of total expression value
#   K  Receptor-expressing cells
#   L  Receptor detection rate
#   M  Receptor average expression value
#   N  Receptor total expression value
#   O  Receptor derived specificity of average expression value
#   P  Receptor derived specificity of total expression value
#   Q  Edge average expression weight
#   R  Edge total expression weight
#   S  Edge average expression derived specificity
#   T  Edge total expression derived specificity

$rows = @(
    # Sending=ECs, Ligand=a, Receptor=Mc5r, Target=ECs
    @("ECs","a","Mc5r","ECs",1,0.3333333333333333,0.045339,0.136017,0.1740293637846656,0.1740293637846656,3,1,1.677453,5.032359,0.2543757648546719,0.2543757648546719,0.07605404156699999,0.684486374103,0.04426885251989626,0.04426885251989626),
    # Sending=ECs, Ligand=a, Receptor=Mc5r, Target=FAPs
    @("ECs","a","Mc5r","FAPs",1,0.3333333333333333,0.045339,0.136017,0.1740293637846656,0.1740293637846656,3,1,3.700487666666666,11.101463,0.5611569328879042,0.5611569328879042,0.167776410319,1.509987692871,0.09765778401383625,0.09765778401383626),
    # Sending=ECs, Ligand=a, Receptor=Mc5r, Target=MuSCs
    @("ECs","a","Mc5r","MuSCs",1,0.3333333333333333,0.045339,0.136017,0.1740293637846656,0.1740293637846656,3,1,1.148238333333333,3.444715,0.174123510034034,0.174123510034034,0.05205997779499999,0.4685398001549999,0.03030260367117577,0.03030260367117577),
    # Sending=ECs, Ligand=a, Receptor=Mc5r, Target=Resolving-Mac
    @("ECs","a","Mc5r","Resolving-Mac",1,0.3333333333333333,0.045339,0.136017,0.1740293637846656,0.1740293637846656,1,0.3333333333333333,0.06821100000000001,0.204633,0.01034379222338988,0.01034379222338988,0.003092618529,0.027833566761,0.001800123579757313,0.001800123579757313),
    # Sending=MuSCs, Ligand=a, Receptor=Mc5r, Target=ECs
    @("MuSCs","a","Mc5r","ECs",1,0.3333333333333333,0.215186,0.645558,0.8259706362153344,0.8259706362153345,3,1,1.677453,5.032359,0.2543757648546719,0.2543757648546719,0.3609644012579999,3.248679611321999,0.2101069123347757,0.2101069123347757),
    # Sending=MuSCs, Ligand=a, Receptor=Mc5r, Target=FAPs
    @("MuSCs","a","Mc5r","FAPs",1,0.3333333333333333,0.215186,0.645558,0.8259706362153344,0.8259706362153345,3,1,3.700487666666666,11.101463,0.5611569328879042,0.5611569328879042,0.7962931390393332,7.166638251353999,0.463499148874068,0.4634991488740681),
    # Sending=MuSCs, Ligand=a, Receptor=Mc5r, Target=MuSCs
    @("MuSCs","a","Mc5r","MuSCs",1,0.3333333333333333,0.215186,0.645558,0.8259706362153344,0.8259706362153345,3,1,1.148238333333333,3.444715,0.174123510034034,0.174123510034034,0.2470848139966666,2.223763325969999,0.1438209063628582,0.1438209063628582),
    # Sending=MuSCs, Ligand=a, Receptor=Mc5r, Target=Resolving-Mac
    @("MuSCs","a","Mc5r","Resolving-Mac",1,0.3333333333333333,0.215186,0.645558,0.8259706362153344,0.8259706362153345,1,0.3333333333333333,0.06821100000000001,0.204633,0.01034379222338988,0.01034379222338988,0.014678052246,0.132102470214,0.008543668643632572,0.00854366864363257)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
